$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Test")

# ---------------------------------------------------------------------------
# Phase 0 - establish shared-string order (16..22) and cell-style order
# (3 -> 4 -> 5) up front, matching the order the original author created
# them in, before filling in the rest of the (style/string-reusing) cells.
# ---------------------------------------------------------------------------

# sstr16
$ws.Range("B23").Value = "One Of Each Phase"
$ws.Range("B23").Font.Bold = $true

# sstr17, style3 (text format)
$ws.Range("B25").Value = "- Growth Phase with known Duration and Growth Rate Factor"
$ws.Range("B25").NumberFormat = "@"

# sstr18, sstr19
$ws.Range("D34").Value = "Time"
$ws.Range("C34").Value = "#"
$ws.Range("E34").Value = "Effect"

# sstr20
$ws.Range("B32").Value = "StepSize"

# sstr21
$ws.Range("B47").Value = "- Steady Phase with known Duration "
$ws.Range("B47").NumberFormat = "@"

# sstr22
$ws.Range("B58").Value = "- Decay Phase with known Duration and Growth Rate Factor"
$ws.Range("B58").NumberFormat = "@"

# style4 (0.00, no alignment)
$ws.Range("C50").Formula = "=E45"
$ws.Range("C50").NumberFormat = "0.00"

# style5 (0.00, centered)
$ws.Range("E35").Formula = "=C28"
$ws.Range("E35").NumberFormat = "0.00"
$ws.Range("E35").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Phase 1 - Growth Phase (rows 23-45)
# ---------------------------------------------------------------------------
$ws.Range("B27").Value = "InitialXt"
$ws.Range("C27").Value = 0
$ws.Range("B28").Value = "InitialYq"
$ws.Range("C28").Value = 0

$ws.Range("B30").Value = "Duration "
$ws.Range("C30").Value = 100
$ws.Range("D30").Value = "s"

$ws.Range("B31").Value = "Growth Rate Factor"
$ws.Range("C31").Value = 0.047
$ws.Range("D31").Value = "kW/s²"

$ws.Range("C32").Value = 10

$ws.Range("C33:E33").HorizontalAlignment = -4108
$ws.Range("C34:E34").HorizontalAlignment = -4108

$ws.Range("C35").Value = 0
$ws.Range("D35").Formula = "=C27"
$ws.Range("C35:D35").HorizontalAlignment = -4108

$ws.Range("C36").Value = 1
$ws.Range("D36").Formula = "=`$C`$32*C36"
$ws.Range("E36").Formula = "=`$C`$31*D36^2"
$ws.Range("C36:D36").HorizontalAlignment = -4108
$ws.Range("E36").NumberFormat = "0.00"
$ws.Range("E36").HorizontalAlignment = -4108

$ws.Range("D37:D45").Formula = "=`$C`$32*C37"
$ws.Range("E37:E45").Formula = "=`$C`$31*D37^2"
$ws.Range("C37:D45").HorizontalAlignment = -4108
$ws.Range("E37:E45").NumberFormat = "0.00"
$ws.Range("E37:E45").HorizontalAlignment = -4108
for ($i = 2; $i -le 10; $i++) {
    $r = 35 + $i
    $ws.Cells.Item($r, 3).Value = $i
}

# ---------------------------------------------------------------------------
# Phase 2 - Steady Phase (rows 47-56)
# ---------------------------------------------------------------------------
$ws.Range("B49").Value = "InitialXt"
$ws.Range("C49").Formula = "=D45"

$ws.Range("B50").Value = "InitialYq"

$ws.Range("B52").Value = "Duration "
$ws.Range("C52").Value = 100
$ws.Range("D52").Value = "s"

$ws.Range("C54").Value = "#"
$ws.Range("D54").Value = "Time"
$ws.Range("E54").Value = "Effect"
$ws.Range("C54:E54").HorizontalAlignment = -4108

$ws.Range("C55").Value = 0
$ws.Range("D55").Formula = "=C49"
$ws.Range("E55").Formula = "=C50"
$ws.Range("C55:D55").HorizontalAlignment = -4108
$ws.Range("E55").NumberFormat = "0.00"
$ws.Range("E55").HorizontalAlignment = -4108

$ws.Range("C56").Value = 1
$ws.Range("D56").Formula = "=D55+C52"
$ws.Range("E56").Formula = "=E55"
$ws.Range("C56:D56").HorizontalAlignment = -4108
$ws.Range("E56").NumberFormat = "0.00"
$ws.Range("E56").HorizontalAlignment = -4108

# ---------------------------------------------------------------------------
# Phase 3 - Decay Phase (rows 58-77)
# ---------------------------------------------------------------------------
$ws.Range("B60").Value = "InitialXt"
$ws.Range("C60").Formula = "=D56"

$ws.Range("B61").Value = "InitialYq"
$ws.Range("C61").Formula = "=E56"
$ws.Range("C61").NumberFormat = "0.00"

$ws.Range("B63").Value = "Duration "
$ws.Range("C63").Value = 100
$ws.Range("D63").Value = "s"

$ws.Range("B64").Value = "Growth Rate Factor"
$ws.Range("C64").Value = 0.047
$ws.Range("D64").Value = "kW/s²"

$ws.Range("C66").Value = "#"
$ws.Range("D66").Value = "Time"
$ws.Range("E66").Value = "Effect"
$ws.Range("C66:E66").HorizontalAlignment = -4108

$ws.Range("C67").Value = 0
$ws.Range("D67").Formula = "=C60"
$ws.Range("E67").Formula = "=C61"
$ws.Range("C67:D67").HorizontalAlignment = -4108
$ws.Range("E67").NumberFormat = "0.00"
$ws.Range("E67").HorizontalAlignment = -4108

$ws.Range("C68").Value = 1
$ws.Range("C68").HorizontalAlignment = -4108
$ws.Range("D68").HorizontalAlignment = -4108
$ws.Range("E68").NumberFormat = "0.00"
$ws.Range("E68").HorizontalAlignment = -4108

$ws.Range("C69:C77").HorizontalAlignment = -4108
$ws.Range("D69:D77").HorizontalAlignment = -4108
$ws.Range("E69:E77").HorizontalAlignment = -4108
for ($i = 2; $i -le 10; $i++) {
    $r = 67 + $i
    $ws.Cells.Item($r, 3).Value = $i
}

# ---------------------------------------------------------------------------
# Sheet view: scrolled/selected state
# ---------------------------------------------------------------------------
$ws.Range("F77:F78").Select()
$excel.ActiveWindow.ScrollRow = 37
